$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The three hyperlinks below the row being removed (all pointing at the same
# "William Ellard" picture URL) need to be re-pointed one row higher, since
# EntireRow.Delete does not itself re-anchor the worksheet's hyperlink list
# in this engine. Capture the shared target URL before editing.
$target = "https://swimming.box.com/shared/static/12ruhuduew34hnxgnbvvb8lq46u4a5y3.png"

$ws.Cells.Item(76, 5).Hyperlinks.Delete()
$ws.Cells.Item(77, 5).Hyperlinks.Delete()
$ws.Cells.Item(78, 5).Hyperlinks.Delete()

# Delete entire row 72 (Toni Shaw / F 100 Fr S9), shifting subsequent rows up.
$ws.Rows.Item(72).Delete()

# Re-create the hyperlinks one row higher than they were.
$ws.Hyperlinks.Add($ws.Cells.Item(75, 5), $target)
$ws.Hyperlinks.Add($ws.Cells.Item(76, 5), $target)
$ws.Hyperlinks.Add($ws.Cells.Item(77, 5), $target)
